$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDGAUGRHC%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T150003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCIGNWYjbzaXNtig5nhtcNnu1%2FI2zGwM0VZSLGyyHyO5YqAiANPoDzzPrWQuFslYQDrIT6UH%2FKNXt60tPxhgfOtZTTwCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMSYfaiDuy7N3JIsZDKroD2dNaveK1rWh3QlRhjIvgvwN9Jp6HV6lEKHZ1sZR%2B1uv7bQEQ9jj1mZecJu5QwpwiI1pl14XryUao9GnlXCIz%2BNwDfryF7WIqiGK9Ig%2BkrakuxoFdJ0qdeY1F7hDLDLzFn3iHg%2Bd8Hamu1c0eTrH3%2Bzy8sH71mOiL6UYXWLYLBNrcxObu%2B0GCXpB%2BrSMgsfsN6dQIyhJrRmUbw2Vw346HfdcurrjqMjAv%2FSSeNd58lVhEgAFzjL2X9lH24tyoatXE0XS8xQ0bDRtwwsYfr0aTQAB8R%2FpKn9%2Fd5Jneeumk%2BTCo31qaldTcCk682hxZcVibIRq9r16Uud0vRMyjy6YqBBOznDrfEmab1qevdSY%2FSvZOzsIkzgcNlqTPYkqYvDb67yvGBB8soLqBibonBshOxOOUcj95k8xcwy9cPdUie7UQrwHLk2MvZBSg%2FAzo1p3XWHoywEiH0ml%2Ft8XN5nNytMJ%2BUz76gmTWrnd8jZWEMifkXj7rEZAVE075qM3%2FTx0bv365ShU%2BSS8xMpRWp9u%2BYpTEnBiEjxR%2BZlrw6KtgqHY6x7gxQ4XLNEmsv%2BRbUMspaiqGlOu%2BZ3xeZTDRtaS%2FBjqmAVP%2FbOnfVwxLUVNPJFPEHNKmdDf9BbAfAAPmEwumBRttANoa2txRdncDEmMeSE8eM9LfUY21SKVrfHm5C34Vwa04k9THyXrzx5%2BMuT63SCl4xljypaMzksED2Apbv9xH2MMLxCtA1Ni9jUudwtLwb%2BEaR23KjzPOtDooABADZY%2FF2v%2FA0MDlDAQiKRLuyYUQS6sfMQBHHFy6I0%2BhTAr%2FvlznpnGYF7I%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2023%3A00%3A03%20GMT&X-Amz-Signature=22d24b16a2d4f3a8f7caa7396cdfbc5a194076bc20c1f7ca9064b1b4115e94ac'

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDGAUGRHC%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T150003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCIGNWYjbzaXNtig5nhtcNnu1%2FI2zGwM0VZSLGyyHyO5YqAiANPoDzzPrWQuFslYQDrIT6UH%2FKNXt60tPxhgfOtZTTwCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMSYfaiDuy7N3JIsZDKroD2dNaveK1rWh3QlRhjIvgvwN9Jp6HV6lEKHZ1sZR%2B1uv7bQEQ9jj1mZecJu5QwpwiI1pl14XryUao9GnlXCIz%2BNwDfryF7WIqiGK9Ig%2BkrakuxoFdJ0qdeY1F7hDLDLzFn3iHg%2Bd8Hamu1c0eTrH3%2Bzy8sH71mOiL6UYXWLYLBNrcxObu%2B0GCXpB%2BrSMgsfsN6dQIyhJrRmUbw2Vw346HfdcurrjqMjAv%2FSSeNd58lVhEgAFzjL2X9lH24tyoatXE0XS8xQ0bDRtwwsYfr0aTQAB8R%2FpKn9%2Fd5Jneeumk%2BTCo31qaldTcCk682hxZcVibIRq9r16Uud0vRMyjy6YqBBOznDrfEmab1qevdSY%2FSvZOzsIkzgcNlqTPYkqYvDb67yvGBB8soLqBibonBshOxOOUcj95k8xcwy9cPdUie7UQrwHLk2MvZBSg%2FAzo1p3XWHoywEiH0ml%2Ft8XN5nNytMJ%2BUz76gmTWrnd8jZWEMifkXj7rEZAVE075qM3%2FTx0bv365ShU%2BSS8xMpRWp9u%2BYpTEnBiEjxR%2BZlrw6KtgqHY6x7gxQ4XLNEmsv%2BRbUMspaiqGlOu%2BZ3xeZTDRtaS%2FBjqmAVP%2FbOnfVwxLUVNPJFPEHNKmdDf9BbAfAAPmEwumBRttANoa2txRdncDEmMeSE8eM9LfUY21SKVrfHm5C34Vwa04k9THyXrzx5%2BMuT63SCl4xljypaMzksED2Apbv9xH2MMLxCtA1Ni9jUudwtLwb%2BEaR23KjzPOtDooABADZY%2FF2v%2FA0MDlDAQiKRLuyYUQS6sfMQBHHFy6I0%2BhTAr%2FvlznpnGYF7I%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2023%3A00%3A03%20GMT&X-Amz-Signature=1f521244548d9b50619821281c35c8b7d01d41edf7d817477af9dacd0ba79d9f'

$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDGAUGRHC%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T150003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCIGNWYjbzaXNtig5nhtcNnu1%2FI2zGwM0VZSLGyyHyO5YqAiANPoDzzPrWQuFslYQDrIT6UH%2FKNXt60tPxhgfOtZTTwCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMSYfaiDuy7N3JIsZDKroD2dNaveK1rWh3QlRhjIvgvwN9Jp6HV6lEKHZ1sZR%2B1uv7bQEQ9jj1mZecJu5QwpwiI1pl14XryUao9GnlXCIz%2BNwDfryF7WIqiGK9Ig%2BkrakuxoFdJ0qdeY1F7hDLDLzFn3iHg%2Bd8Hamu1c0eTrH3%2Bzy8sH71mOiL6UYXWLYLBNrcxObu%2B0GCXpB%2BrSMgsfsN6dQIyhJrRmUbw2Vw346HfdcurrjqMjAv%2FSSeNd58lVhEgAFzjL2X9lH24tyoatXE0XS8xQ0bDRtwwsYfr0aTQAB8R%2FpKn9%2Fd5Jneeumk%2BTCo31qaldTcCk682hxZcVibIRq9r16Uud0vRMyjy6YqBBOznDrfEmab1qevdSY%2FSvZOzsIkzgcNlqTPYkqYvDb67yvGBB8soLqBibonBshOxOOUcj95k8xcwy9cPdUie7UQrwHLk2MvZBSg%2FAzo1p3XWHoywEiH0ml%2Ft8XN5nNytMJ%2BUz76gmTWrnd8jZWEMifkXj7rEZAVE075qM3%2FTx0bv365ShU%2BSS8xMpRWp9u%2BYpTEnBiEjxR%2BZlrw6KtgqHY6x7gxQ4XLNEmsv%2BRbUMspaiqGlOu%2BZ3xeZTDRtaS%2FBjqmAVP%2FbOnfVwxLUVNPJFPEHNKmdDf9BbAfAAPmEwumBRttANoa2txRdncDEmMeSE8eM9LfUY21SKVrfHm5C34Vwa04k9THyXrzx5%2BMuT63SCl4xljypaMzksED2Apbv9xH2MMLxCtA1Ni9jUudwtLwb%2BEaR23KjzPOtDooABADZY%2FF2v%2FA0MDlDAQiKRLuyYUQS6sfMQBHHFy6I0%2BhTAr%2FvlznpnGYF7I%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2023%3A00%3A03%20GMT&X-Amz-Signature=3a05a21fa4245210610aba4a23079fa712b4cc07fdf14a6bd755b32c929f4dfb'

$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDGAUGRHC%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T150003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCIGNWYjbzaXNtig5nhtcNnu1%2FI2zGwM0VZSLGyyHyO5YqAiANPoDzzPrWQuFslYQDrIT6UH%2FKNXt60tPxhgfOtZTTwCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMSYfaiDuy7N3JIsZDKroD2dNaveK1rWh3QlRhjIvgvwN9Jp6HV6lEKHZ1sZR%2B1uv7bQEQ9jj1mZecJu5QwpwiI1pl14XryUao9GnlXCIz%2BNwDfryF7WIqiGK9Ig%2BkrakuxoFdJ0qdeY1F7hDLDLzFn3iHg%2Bd8Hamu1c0eTrH3%2Bzy8sH71mOiL6UYXWLYLBNrcxObu%2B0GCXpB%2BrSMgsfsN6dQIyhJrRmUbw2Vw346HfdcurrjqMjAv%2FSSeNd58lVhEgAFzjL2X9lH24tyoatXE0XS8xQ0bDRtwwsYfr0aTQAB8R%2FpKn9%2Fd5Jneeumk%2BTCo31qaldTcCk682hxZcVibIRq9r16Uud0vRMyjy6YqBBOznDrfEmab1qevdSY%2FSvZOzsIkzgcNlqTPYkqYvDb67yvGBB8soLqBibonBshOxOOUcj95k8xcwy9cPdUie7UQrwHLk2MvZBSg%2FAzo1p3XWHoywEiH0ml%2Ft8XN5nNytMJ%2BUz76gmTWrnd8jZWEMifkXj7rEZAVE075qM3%2FTx0bv365ShU%2BSS8xMpRWp9u%2BYpTEnBiEjxR%2BZlrw6KtgqHY6x7gxQ4XLNEmsv%2BRbUMspaiqGlOu%2BZ3xeZTDRtaS%2FBjqmAVP%2FbOnfVwxLUVNPJFPEHNKmdDf9BbAfAAPmEwumBRttANoa2txRdncDEmMeSE8eM9LfUY21SKVrfHm5C34Vwa04k9THyXrzx5%2BMuT63SCl4xljypaMzksED2Apbv9xH2MMLxCtA1Ni9jUudwtLwb%2BEaR23KjzPOtDooABADZY%2FF2v%2FA0MDlDAQiKRLuyYUQS6sfMQBHHFy6I0%2BhTAr%2FvlznpnGYF7I%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2023%3A00%3A03%20GMT&X-Amz-Signature=3c8859bde447b92a93f59ba275515ed10b45d1cb58e60b5ca89b596b5d890f2b'

$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDGAUGRHC%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T150003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECMaCXVzLXdlc3QtMiJGMEQCIGNWYjbzaXNtig5nhtcNnu1%2FI2zGwM0VZSLGyyHyO5YqAiANPoDzzPrWQuFslYQDrIT6UH%2FKNXt60tPxhgfOtZTTwCrmAwiL%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMSYfaiDuy7N3JIsZDKroD2dNaveK1rWh3QlRhjIvgvwN9Jp6HV6lEKHZ1sZR%2B1uv7bQEQ9jj1mZecJu5QwpwiI1pl14XryUao9GnlXCIz%2BNwDfryF7WIqiGK9Ig%2BkrakuxoFdJ0qdeY1F7hDLDLzFn3iHg%2Bd8Hamu1c0eTrH3%2Bzy8sH71mOiL6UYXWLYLBNrcxObu%2B0GCXpB%2BrSMgsfsN6dQIyhJrRmUbw2Vw346HfdcurrjqMjAv%2FSSeNd58lVhEgAFzjL2X9lH24tyoatXE0XS8xQ0bDRtwwsYfr0aTQAB8R%2FpKn9%2Fd5Jneeumk%2BTCo31qaldTcCk682hxZcVibIRq9r16Uud0vRMyjy6YqBBOznDrfEmab1qevdSY%2FSvZOzsIkzgcNlqTPYkqYvDb67yvGBB8soLqBibonBshOxOOUcj95k8xcwy9cPdUie7UQrwHLk2MvZBSg%2FAzo1p3XWHoywEiH0ml%2Ft8XN5nNytMJ%2BUz76gmTWrnd8jZWEMifkXj7rEZAVE075qM3%2FTx0bv365ShU%2BSS8xMpRWp9u%2BYpTEnBiEjxR%2BZlrw6KtgqHY6x7gxQ4XLNEmsv%2BRbUMspaiqGlOu%2BZ3xeZTDRtaS%2FBjqmAVP%2FbOnfVwxLUVNPJFPEHNKmdDf9BbAfAAPmEwumBRttANoa2txRdncDEmMeSE8eM9LfUY21SKVrfHm5C34Vwa04k9THyXrzx5%2BMuT63SCl4xljypaMzksED2Apbv9xH2MMLxCtA1Ni9jUudwtLwb%2BEaR23KjzPOtDooABADZY%2FF2v%2FA0MDlDAQiKRLuyYUQS6sfMQBHHFy6I0%2BhTAr%2FvlznpnGYF7I%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2023%3A00%3A03%20GMT&X-Amz-Signature=1194de5a916c3a168f935e9f9fc06858477dfc33574f5e22d4207d3738150be1'
